$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.47568978777511184
$ws.Range("B1").Value = 0.21298007248749995
$ws.Range("A2").Value = 0.3303780960716638
$ws.Range("B2").Value = 0.24623383579351343
